$wb = $excel.ActiveWorkbook

# --- "suggest" sheet: remove "Chao co" / "Sinh hoat chu nhiem" rows for each grade ---
$suggest = $wb.Worksheets.Item("suggest")

# Rows are, before any deletion: 2,3 (grade10) / 18,19 (grade11) / 34,35 (grade12).
# Delete from the bottom up so earlier row numbers stay valid.
$suggest.Rows.Item(35).Delete()
$suggest.Rows.Item(34).Delete()
$suggest.Rows.Item(19).Delete()
$suggest.Rows.Item(18).Delete()
$suggest.Rows.Item(3).Delete()
$suggest.Rows.Item(2).Delete()

# --- make "suggest" the active sheet/tab, matching the new selection state ---
$suggest.Activate()
$suggest.Range("F32").Select()

# "teacher" sheet no longer keeps the tab-selected flag; its own selection stays as-is.
$teacher = $wb.Worksheets.Item("teacher")
$teacher.Range("F3").Select()

# restore "suggest" as the active sheet after touching "teacher"
$suggest.Activate()
$suggest.Range("F32").Select()
